# Add a new "UK" test-data worksheet, modeled on the existing "Poland" sheet,
# and make it the active sheet/selection (matching the authored diff).

$wb = $excel.ActiveWorkbook

# The "Poland" sheet is the template for every per-market sheet in this
# workbook (same layout/styles/merges) - copy it and drop the copy right
# after it, then rename to "UK".
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland)
$uk = $wb.Worksheets.Item($poland.Index + 1)
$uk.Name = "UK"

# Update the market-specific cells. Set B4 (the NGC ticket reference) before
# B2 (the market name) so new shared-string entries are appended in the same
# order as the source edit: "NGC-2741/T3363/T3354" then "UK Market".
$uk.Range("B4").Value = "NGC-2741/T3363/T3354"
$uk.Range("B2").Value = "UK Market"

# Make the new sheet the active tab with B2 selected, like the author left it.
$uk.Activate()
$uk.Range("B2").Select()
